$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("N6").ClearContents()
$ws.Range("H6").Value = 314.7
$ws.Range("I6").Value = 314.7
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 944.0999999999999
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -832.0999999999999
$ws.Range("M76").ClearContents()
$ws.Range("H76").Value = 0
$ws.Range("I76").Value = 0
$ws.Range("K76").Value = 0
$ws.Range("M79").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("I79").Value = 0
$ws.Range("K79").Value = 0
$ws.Range("H125").Value = 1000
$ws.Range("I125").Value = 1000
$ws.Range("K125").Value = 9000
$ws.Range("M125").Value = -6540
$ws.Range("H138").Value = 1986.75
$ws.Range("I138").Value = 683.625
$ws.Range("J138").Value = 3289.875
$ws.Range("K138").Value = 2050.875
$ws.Range("L138").Value = 9869.625
$ws.Range("M138").Value = 3089.125
$ws.Range("N138").Value = -20149.625

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 829
$ws.Range("I110").Value = 829
$ws.Range("K110").Value = 829
$ws.Range("M110").Value = 1216
$ws.Range("H119").Value = 44500
$ws.Range("J119").Value = 44500
$ws.Range("L119").Value = 44500
$ws.Range("N119").Value = -54176
$ws.Range("N122").ClearContents()
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 3000
$ws.Range("I105").Value = 3000
$ws.Range("K105").Value = 3000
$ws.Range("M105").Value = -1253

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("N5").ClearContents()
$ws.Range("H5").Value = 578
$ws.Range("I5").Value = 578
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 578
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -466
$ws.Range("H39").Value = 27033.334
$ws.Range("I39").Value = 1100
$ws.Range("K39").Value = 1100
$ws.Range("M39").Value = -709
$ws.Range("H49").Value = 27033.334
$ws.Range("I49").Value = 1100
$ws.Range("K49").Value = 1100
$ws.Range("M49").Value = -918
$ws.Range("H86").Value = 8800.625
$ws.Range("I86").Value = 9327.5
$ws.Range("K86").Value = 9327.5
$ws.Range("M86").Value = -8204.5
$ws.Range("H89").Value = 8800.625
$ws.Range("I89").Value = 9327.5
$ws.Range("K89").Value = 46637.5
$ws.Range("M89").Value = -41021.5
$ws.Range("H132").Value = 1366.3334
$ws.Range("I132").Value = 1799.5
$ws.Range("J132").Value = 500
$ws.Range("K132").Value = 5398.5
$ws.Range("L132").Value = 1500
$ws.Range("M132").Value = -2868.5
$ws.Range("N132").Value = -6560
$ws.Range("H134").Value = 1994
$ws.Range("I134").Value = 1994
$ws.Range("K134").Value = 5982
$ws.Range("M134").Value = -3447

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 12.62963
$ws.Range("I2").Value = 7.368421
$ws.Range("J2").Value = 25.125
$ws.Range("K2").Value = 44.210526
$ws.Range("L2").Value = 150.75
$ws.Range("M2").Value = 68.789474
$ws.Range("N2").Value = -376.75
$ws.Range("H3").Value = 750
$ws.Range("I3").Value = 1000
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = -2888
$ws.Range("N3").Value = -1724
$ws.Range("H17").Value = 431.66666
$ws.Range("I17").Value = 240
$ws.Range("K17").Value = 720
$ws.Range("M17").Value = -551
$ws.Range("H34").Value = 568.2222
$ws.Range("J34").Value = 2225
$ws.Range("L34").Value = 6675
$ws.Range("N34").Value = -6843
$ws.Range("H39").Value = 22000
$ws.Range("J39").Value = 27000
$ws.Range("L39").Value = 81000
$ws.Range("N39").Value = -81588
$ws.Range("H55").Value = 1900
$ws.Range("J55").Value = 3550
$ws.Range("L55").Value = 10650
$ws.Range("N55").Value = -11004

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3927
$ws.Range("I70").Value = 3933
$ws.Range("K70").Value = 3933
$ws.Range("M70").Value = -3663
$ws.Range("H73").Value = 3927
$ws.Range("I73").Value = 3933
$ws.Range("K73").Value = 3933
$ws.Range("M73").Value = -2997

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3502.5
$ws.Range("I7").Value = 5000
$ws.Range("J7").Value = 2005
$ws.Range("K7").Value = 5000
$ws.Range("L7").Value = 2005
$ws.Range("M7").Value = -4888
$ws.Range("N7").Value = -2229
$ws.Range("M21").ClearContents()
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("H22").Value = 3494.4666
$ws.Range("I22").Value = 2574.75
$ws.Range("J22").Value = 3828.9092
$ws.Range("K22").Value = 2574.75
$ws.Range("L22").Value = 3828.9092
$ws.Range("M22").Value = -2279.75
$ws.Range("N22").Value = -4418.9092
$ws.Range("H24").Value = 22500
$ws.Range("I24").Value = 15000
$ws.Range("J24").Value = 30000
$ws.Range("K24").Value = 15000
$ws.Range("L24").Value = 30000
$ws.Range("M24").Value = -14657
$ws.Range("N24").Value = -30686
$ws.Range("H27").Value = 3494.4666
$ws.Range("I27").Value = 2574.75
$ws.Range("J27").Value = 3828.9092
$ws.Range("K27").Value = 2574.75
$ws.Range("L27").Value = 3828.9092
$ws.Range("M27").Value = -2467.75
$ws.Range("N27").Value = -4042.9092
$ws.Range("H46").Value = 255123.75
$ws.Range("J46").Value = 5855.7144
$ws.Range("L46").Value = 5855.7144
$ws.Range("N46").Value = -6231.7144
$ws.Range("H61").Value = 2365.5715
$ws.Range("I61").Value = 3030.8
$ws.Range("J61").Value = 702.5
$ws.Range("K61").Value = 3030.8
$ws.Range("L61").Value = 702.5
$ws.Range("M61").Value = -2828.8
$ws.Range("N61").Value = -1106.5
$ws.Range("H113").Value = 2365.5715
$ws.Range("I113").Value = 3030.8
$ws.Range("J113").Value = 702.5
$ws.Range("K113").Value = 3030.8
$ws.Range("L113").Value = 702.5
$ws.Range("M113").Value = -860.8000000000002
$ws.Range("N113").Value = -5042.5
$ws.Range("H126").Value = 3502.5
$ws.Range("I126").Value = 5000
$ws.Range("J126").Value = 2005
$ws.Range("K126").Value = 15000
$ws.Range("L126").Value = 6015
$ws.Range("M126").Value = -12530
$ws.Range("N126").Value = -10955
$ws.Range("H136").Value = 4723
$ws.Range("I136").Value = 5000
$ws.Range("J136").Value = 4584.5
$ws.Range("K136").Value = 15000
$ws.Range("L136").Value = 13753.5
$ws.Range("M136").Value = -12450
$ws.Range("N136").Value = -18853.5
$ws.Range("H137").Value = 49995
$ws.Range("I137").Value = 49995
$ws.Range("K137").Value = 49995
$ws.Range("M137").Value = -44895

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("M12").ClearContents()
$ws.Range("H12").Value = 200
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 200
$ws.Range("N12").Value = -484
$ws.Range("H132").Value = 1438.8572
$ws.Range("I132").Value = 1438.8572
$ws.Range("K132").Value = 4316.571599999999
$ws.Range("M132").Value = -1786.571599999999
$ws.Range("H136").Value = 20000
$ws.Range("I136").Value = 20000
$ws.Range("K136").Value = 60000
$ws.Range("M136").Value = -57450
